# Update "想去人数" (want-to-go count) figures across the three sheets
# that track this workbook's exhibition/performance listings.
#
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 764
$ws1.Range("F5").Value = 181
$ws1.Range("F6").Value = 43
$ws1.Range("F7").Value = 1067
$ws1.Range("F8").Value = 643
$ws1.Range("F9").Value = 735
$ws1.Range("F10").Value = 1281
$ws1.Range("F11").Value = 258
$ws1.Range("F12").Value = 992
$ws1.Range("F13").Value = 48
$ws1.Range("F19").Value = 518
$ws1.Range("F21").Value = 729
$ws1.Range("F22").Value = 211
$ws1.Range("F23").Value = 146

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 961
$ws2.Range("F5").Value = 203

# Sheet "全部类型" (All types - combined listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 764
$ws4.Range("F7").Value = 181
$ws4.Range("F8").Value = 961
$ws4.Range("F9").Value = 43
$ws4.Range("F10").Value = 1067
$ws4.Range("F11").Value = 643
$ws4.Range("F12").Value = 735
$ws4.Range("F13").Value = 1281
$ws4.Range("F14").Value = 258
$ws4.Range("F15").Value = 992
$ws4.Range("F16").Value = 48
$ws4.Range("F21").Value = 203
$ws4.Range("F27").Value = 518
$ws4.Range("F29").Value = 729
$ws4.Range("F30").Value = 211
$ws4.Range("F32").Value = 146
